$d = $word.ActiveDocument

# Helper characters for the curly quotes used throughout the document.
$lq = [char]0x201C   # “
$rq = [char]0x201D   # ”
$ap = [char]0x2019   # ’

# --- 1. "1. " step becomes its own sentence, and a brand-new step "2. System
#        displays an empty form" is inserted right after it (all the other
#        numbered steps shift down by one). ------------------------------
$d.Content.Find.Execute(
    "1. System displays an empty form", $true, $false, $false, $false, $false,
    $true, 1, $false, "1. The user chooses to create a movie", 2) | Out-Null

# Find the paragraph we just edited and insert a new paragraph after it that
# carries the text the "1." step used to have, renumbered to "2.".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd([char]13,[char]7) -eq "1. The user chooses to create a movie") {
        $para.Range.InsertParagraphAfter() | Out-Null
        $newPara = $d.Paragraphs($i + 1)
        $newPara.Alignment = 3
        $newPara.Range.Text = "2. System displays an empty form"
        break
    }
}

# --- 2. Renumber the rest of the Main Success Scenario steps. ------------
$d.Content.Find.Execute(
    "2. User fills in fields with data", $true, $false, $false, $false, $false,
    $true, 1, $false, "3. User fills in fields with data", 2) | Out-Null

$d.Content.Find.Execute(
    "3. Movie is created with specified data", $true, $false, $false, $false, $false,
    $true, 1, $false, "4. Movie is created with specified data", 2) | Out-Null

$d.Content.Find.Execute(
    "4. Movie is added to a list", $true, $false, $false, $false, $false,
    $true, 1, $false, "5. Movie is added to a list", 2) | Out-Null

# --- 3. Renumber the Extensions entries that referenced step "2". --------
$d.Content.Find.Execute(
    "2a. The user inputs the wrong format of the year of creation", $true, $false, $false, $false, $false,
    $true, 1, $false, "3a. The user inputs the wrong format of the year of creation", 2) | Out-Null

$d.Content.Find.Execute(
    "2b The user inputs a year of creation that is further than the current year", $true, $false, $false, $false, $false,
    $true, 1, $false, "3b The user inputs a year of creation that is further than the current year", 2) | Out-Null

$d.Content.Find.Execute(
    "2c The user inputs the wrong format of the release date", $true, $false, $false, $false, $false,
    $true, 1, $false, "3c The user inputs the wrong format of the release date", 2) | Out-Null

$d.Content.Find.Execute(
    "2d The user inputs the release date", $true, $false, $false, $false, $false,
    $true, 1, $false, "3d The user inputs the release date", 2) | Out-Null

$d.Content.Find.Execute(
    "2f The  user inputs the wrong format of the price", $true, $false, $false, $false, $false,
    $true, 1, $false, "3f The  user inputs the wrong format of the price", 2) | Out-Null

# --- 4. Renumber the "Movie with the same X" extensions from 3a/3a to 4a/4b.
$d.Content.Find.Execute(
    "3a.Movie with the same title is already created", $true, $false, $false, $false, $false,
    $true, 1, $false, "4a.Movie with the same title is already created", 2) | Out-Null

$d.Content.Find.Execute(
    "3a.Movie with the same description is already created", $true, $false, $false, $false, $false,
    $true, 1, $false, "4b.Movie with the same description is already created", 2) | Out-Null
